$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scheduled-task refresh: row 12's timestamp gets a higher-precision
# re-stamp, and a brand-new row 13 is appended with the next reading.

$ws.Cells.Item(12, 1).Value = 45865.54202207176

$ws.Cells.Item(13, 1).Value = 45865.5836164567
$ws.Cells.Item(13, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 2).Value = 2025
$ws.Cells.Item(13, 3).Value = 30
$ws.Cells.Item(13, 4).Value = 18.06
$ws.Cells.Item(13, 5).Value = 77.70999999999999
$ws.Cells.Item(13, 6).Value = 119.53
$ws.Cells.Item(13, 7).Value = 12.35
$ws.Cells.Item(13, 8).Value = "ESE"
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = "14:00:24"
